$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.18739966666666
$ws.Range("H2").Value = 189.562199
$ws.Range("I2").Value = 0.09596345243430386
$ws.Range("J2").Value = 0.09988075390087989
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1744923333333333
$ws.Range("N2").Value = 0.523477
$ws.Range("O2").Value = 0.07822917822503123
$ws.Range("P2").Value = 0.08239975633156223
$ws.Range("Q2").Value = 11.02571680510255
$ws.Range("R2").Value = 99.231451245923
$ws.Range("S2").Value = 0.007507142023572463
$ws.Range("T2").Value = 0.008230149783645236
$ws.Range("G3").Value = 63.18739966666666
$ws.Range("H3").Value = 189.562199
$ws.Range("I3").Value = 0.09596345243430386
$ws.Range("J3").Value = 0.09988075390087989
$ws.Range("O3").Value = 0.5425629886502931
$ws.Range("P3").Value = 0.5714882742434749
$ws.Range("Q3").Value = 76.46949638893277
$ws.Range("R3").Value = 688.225467500395
$ws.Range("S3").Value = 0.05206621755395614
$ws.Range("T3").Value = 0.05708067967695107
$ws.Range("G4").Value = 63.18739966666666
$ws.Range("H4").Value = 189.562199
$ws.Range("I4").Value = 0.09596345243430386
$ws.Range("J4").Value = 0.09988075390087989
$ws.Range("M4").Value = 0.1427166666666667
$ws.Range("N4").Value = 0.42815
$ws.Range("O4").Value = 0.06398337015197826
$ws.Range("P4").Value = 0.06739447133944447
$ws.Range("Q4").Value = 9.017895055761112
$ws.Range("R4").Value = 81.16105550185
$ws.Range("S4").Value = 0.006140065098165824
$ws.Range("T4").Value = 0.006731410606134956
$ws.Range("G5").Value = 63.18739966666666
$ws.Range("H5").Value = 189.562199
$ws.Range("I5").Value = 0.09596345243430386
$ws.Range("J5").Value = 0.09988075390087989
$ws.Range("M5").Value = 0.3386875
$ws.Range("N5").Value = 0.6773750000000001
$ws.Range("O5").Value = 0.15184188493529
$ws.Range("P5").Value = 0.1066246175956001
$ws.Range("Q5").Value = 21.40078242460417
$ws.Range("R5").Value = 128.404694547625
$ws.Range("S5").Value = 0.01457127150252274
$ws.Range("T5").Value = 0.01064974718984157
$ws.Range("G6").Value = 63.18739966666666
$ws.Range("H6").Value = 189.562199
$ws.Range("I6").Value = 0.09596345243430386
$ws.Range("J6").Value = 0.09988075390087989
$ws.Range("M6").Value = 0.3644293333333333
$ws.Range("N6").Value = 1.093288
$ws.Range("O6").Value = 0.1633825780374074
$ws.Range("P6").Value = 0.1720928804899184
$ws.Range("Q6").Value = 23.02734193559022
$ws.Range("R6").Value = 207.246077420312
$ws.Range("S6").Value = 0.01567875625608669
$ws.Range("T6").Value = 0.01718876664430707
$ws.Range("I7").Value = 0.3063997713314046
$ws.Range("J7").Value = 0.3189072441572365
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1744923333333333
$ws.Range("N7").Value = 0.523477
$ws.Range("O7").Value = 0.07822917822503123
$ws.Range("P7").Value = 0.08239975633156223
$ws.Range("Q7").Value = 35.20378875656856
$ws.Range("R7").Value = 316.834098809117
$ws.Range("S7").Value = 0.02396940231959327
$ws.Range("T7").Value = 0.02627787921092631
$ws.Range("I8").Value = 0.3063997713314046
$ws.Range("J8").Value = 0.3189072441572365
$ws.Range("O8").Value = 0.5425629886502931
$ws.Range("P8").Value = 0.5714882742434749
$ws.Range("S8").Value = 0.1662411756553333
$ws.Range("T8").Value = 0.1822517506071616
$ws.Range("I9").Value = 0.3063997713314046
$ws.Range("J9").Value = 0.3189072441572365
$ws.Range("M9").Value = 0.1427166666666667
$ws.Range("N9").Value = 0.42815
$ws.Range("O9").Value = 0.06398337015197826
$ws.Range("P9").Value = 0.06739447133944447
$ws.Range("Q9").Value = 28.79305519846111
$ws.Range("R9").Value = 259.13749678615
$ws.Range("S9").Value = 0.01960448998357876
$ws.Range("T9").Value = 0.0214925851262961
$ws.Range("I10").Value = 0.3063997713314046
$ws.Range("J10").Value = 0.3189072441572365
$ws.Range("M10").Value = 0.3386875
$ws.Range("N10").Value = 0.6773750000000001
$ws.Range("O10").Value = 0.15184188493529
$ws.Range("P10").Value = 0.1066246175956001
$ws.Range("Q10").Value = 68.33012646872918
$ws.Range("R10").Value = 409.980758812375
$ws.Range("S10").Value = 0.0465243188227023
$ws.Range("T10").Value = 0.03400336295673204
$ws.Range("I11").Value = 0.3063997713314046
$ws.Range("J11").Value = 0.3189072441572365
$ws.Range("M11").Value = 0.3644293333333333
$ws.Range("N11").Value = 1.093288
$ws.Range("O11").Value = 0.1633825780374074
$ws.Range("P11").Value = 0.1720928804899184
$ws.Range("Q11").Value = 73.52353551749422
$ws.Range("R11").Value = 661.711819657448
$ws.Range("S11").Value = 0.05006038455019701
$ws.Range("T11").Value = 0.05488166625612053
$ws.Range("G12").Value = 170.2928416666667
$ws.Range("H12").Value = 510.878525
$ws.Range("I12").Value = 0.2586257560429799
$ws.Range("J12").Value = 0.2691830570543736
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1744923333333333
$ws.Range("N12").Value = 0.523477
$ws.Range("O12").Value = 0.07822917822503123
$ws.Range("P12").Value = 0.08239975633156223
$ws.Range("Q12").Value = 29.71479529238056
$ws.Range("R12").Value = 267.433157631425
$ws.Range("S12").Value = 0.02023208036306972
$ws.Range("T12").Value = 0.02218061830986539
$ws.Range("G13").Value = 170.2928416666667
$ws.Range("H13").Value = 510.878525
$ws.Range("I13").Value = 0.2586257560429799
$ws.Range("J13").Value = 0.2691830570543736
$ws.Range("O13").Value = 0.5425629886502931
$ws.Range("P13").Value = 0.5714882742434749
$ws.Range("Q13").Value = 206.0886808064028
$ws.Range("R13").Value = 1854.798127257625
$ws.Range("S13").Value = 0.1403207631406208
$ws.Range("T13").Value = 0.1538349607315868
$ws.Range("G14").Value = 170.2928416666667
$ws.Range("H14").Value = 510.878525
$ws.Range("I14").Value = 0.2586257560429799
$ws.Range("J14").Value = 0.2691830570543736
$ws.Range("M14").Value = 0.1427166666666667
$ws.Range("N14").Value = 0.42815
$ws.Range("O14").Value = 0.06398337015197826
$ws.Range("P14").Value = 0.06739447133944447
$ws.Range("Q14").Value = 24.30362671986111
$ws.Range("R14").Value = 218.73264047875
$ws.Range("S14").Value = 0.01654774747973321
$ws.Range("T14").Value = 0.01814144982371502
$ws.Range("G15").Value = 170.2928416666667
$ws.Range("H15").Value = 510.878525
$ws.Range("I15").Value = 0.2586257560429799
$ws.Range("J15").Value = 0.2691830570543736
$ws.Range("M15").Value = 0.3386875
$ws.Range("N15").Value = 0.6773750000000001
$ws.Range("O15").Value = 0.15184188493529
$ws.Range("P15").Value = 0.1066246175956001
$ws.Range("Q15").Value = 57.67605681197917
$ws.Range("R15").Value = 346.056340871875
$ws.Range("S15").Value = 0.03927022229038054
$ws.Range("T15").Value = 0.0287015405216372
$ws.Range("G16").Value = 170.2928416666667
$ws.Range("H16").Value = 510.878525
$ws.Range("I16").Value = 0.2586257560429799
$ws.Range("J16").Value = 0.2691830570543736
$ws.Range("M16").Value = 0.3644293333333333
$ws.Range("N16").Value = 1.093288
$ws.Range("O16").Value = 0.1633825780374074
$ws.Range("P16").Value = 0.1720928804899184
$ws.Range("Q16").Value = 62.05970676002222
$ws.Range("R16").Value = 558.5373608402
$ws.Range("S16").Value = 0.04225494276917567
$ws.Range("T16").Value = 0.04632448766756919
$ws.Range("G17").Value = 77.473122
$ws.Range("H17").Value = 154.946244
$ws.Range("I17").Value = 0.1176593481802354
$ws.Range("J17").Value = 0.08164152846121862
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1744923333333333
$ws.Range("N17").Value = 0.523477
$ws.Range("O17").Value = 0.07822917822503123
$ws.Range("P17").Value = 0.08239975633156223
$ws.Range("Q17").Value = 13.518465828398
$ws.Range("R17").Value = 81.110794970388
$ws.Range("S17").Value = 0.009204394118632634
$ws.Range("T17").Value = 0.006727242051740717
$ws.Range("G18").Value = 77.473122
$ws.Range("H18").Value = 154.946244
$ws.Range("I18").Value = 0.1176593481802354
$ws.Range("J18").Value = 0.08164152846121862
$ws.Range("O18").Value = 0.5425629886502931
$ws.Range("P18").Value = 0.5714882742434749
$ws.Range("Q18").Value = 93.75810136627
$ws.Range("R18").Value = 562.5486081976201
$ws.Range("S18").Value = 0.06383760759131392
$ws.Range("T18").Value = 0.04665717620690137
$ws.Range("G19").Value = 77.473122
$ws.Range("H19").Value = 154.946244
$ws.Range("I19").Value = 0.1176593481802354
$ws.Range("J19").Value = 0.08164152846121862
$ws.Range("M19").Value = 0.1427166666666667
$ws.Range("N19").Value = 0.42815
$ws.Range("O19").Value = 0.06398337015197826
$ws.Range("P19").Value = 0.06739447133944447
$ws.Range("Q19").Value = 11.0567057281
$ws.Range("R19").Value = 66.34023436860001
$ws.Range("S19").Value = 0.007528241626456488
$ws.Range("T19").Value = 0.005502187649988039
$ws.Range("G20").Value = 77.473122
$ws.Range("H20").Value = 154.946244
$ws.Range("I20").Value = 0.1176593481802354
$ws.Range("J20").Value = 0.08164152846121862
$ws.Range("M20").Value = 0.3386875
$ws.Range("N20").Value = 0.6773750000000001
$ws.Range("O20").Value = 0.15184188493529
$ws.Range("P20").Value = 0.1066246175956001
$ws.Range("Q20").Value = 26.239178007375
$ws.Range("R20").Value = 104.9567120295
$ws.Range("S20").Value = 0.01786561720794452
$ws.Range("T20").Value = 0.008704996752097741
$ws.Range("G21").Value = 77.473122
$ws.Range("H21").Value = 154.946244
$ws.Range("I21").Value = 0.1176593481802354
$ws.Range("J21").Value = 0.08164152846121862
$ws.Range("M21").Value = 0.3644293333333333
$ws.Range("N21").Value = 1.093288
$ws.Range("O21").Value = 0.1633825780374074
$ws.Range("P21").Value = 0.1720928804899184
$ws.Range("Q21").Value = 28.233478201712
$ws.Range("R21").Value = 169.400869210272
$ws.Range("S21").Value = 0.01922348763588779
$ws.Range("T21").Value = 0.01404992580049077
$ws.Range("G22").Value = 145.7496183333334
$ws.Range("H22").Value = 437.248855
$ws.Range("I22").Value = 0.2213516720110761
$ws.Range("J22").Value = 0.2303874164262914
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.1744923333333333
$ws.Range("N22").Value = 0.523477
$ws.Range("O22").Value = 0.07822917822503123
$ws.Range("P22").Value = 0.08239975633156223
$ws.Range("Q22").Value = 25.43219098542611
$ws.Range("R22").Value = 228.889718868835
$ws.Range("S22").Value = 0.01731615940016313
$ws.Range("T22").Value = 0.01898386697538457
$ws.Range("G23").Value = 145.7496183333334
$ws.Range("H23").Value = 437.248855
$ws.Range("I23").Value = 0.2213516720110761
$ws.Range("J23").Value = 0.2303874164262914
$ws.Range("O23").Value = 0.5425629886502931
$ws.Range("P23").Value = 0.5714882742434749
$ws.Range("Q23").Value = 176.3864310230306
$ws.Range("R23").Value = 1587.477879207275
$ws.Range("S23").Value = 0.1200972247090689
$ws.Range("T23").Value = 0.1316637070208741
$ws.Range("G24").Value = 145.7496183333334
$ws.Range("H24").Value = 437.248855
$ws.Range("I24").Value = 0.2213516720110761
$ws.Range("J24").Value = 0.2303874164262914
$ws.Range("M24").Value = 0.1427166666666667
$ws.Range("N24").Value = 0.42815
$ws.Range("O24").Value = 0.06398337015197826
$ws.Range("P24").Value = 0.06739447133944447
$ws.Range("Q24").Value = 20.80089969647223
$ws.Range("R24").Value = 187.20809726825
$ws.Range("S24").Value = 0.01416282596404397
$ws.Range("T24").Value = 0.01552683813331035
$ws.Range("G25").Value = 145.7496183333334
$ws.Range("H25").Value = 437.248855
$ws.Range("I25").Value = 0.2213516720110761
$ws.Range("J25").Value = 0.2303874164262914
$ws.Range("M25").Value = 0.3386875
$ws.Range("N25").Value = 0.6773750000000001
$ws.Range("O25").Value = 0.15184188493529
$ws.Range("P25").Value = 0.1066246175956001
$ws.Range("Q25").Value = 49.36357385927084
$ws.Range("R25").Value = 296.1814431556251
$ws.Range("S25").Value = 0.03361045511173986
$ws.Range("T25").Value = 0.02456497017529161
$ws.Range("G26").Value = 145.7496183333334
$ws.Range("H26").Value = 437.248855
$ws.Range("I26").Value = 0.2213516720110761
$ws.Range("J26").Value = 0.2303874164262914
$ws.Range("M26").Value = 0.3644293333333333
$ws.Range("N26").Value = 1.093288
$ws.Range("O26").Value = 0.1633825780374074
$ws.Range("P26").Value = 0.1720928804899184
$ws.Range("Q26").Value = 53.11543624280446
$ws.Range("R26").Value = 478.0389261852401
$ws.Range("S26").Value = 0.03616500682606025
$ws.Range("T26").Value = 0.03964803412143083
